$d = $word.ActiveDocument

# 1. Trim the trailing space from the "Updated figures..." paragraph text.
$d.Content.Find.Execute("to account for this. ", $true, $false, $false, $false, $false, $true, 1, $false, "to account for this.", 2)

# 2. Append two new list-paragraph entries after the last paragraph, matching
#    the formatting (ListParagraph style / numId 7 / Courier New) that
#    InsertParagraphAfter inherits automatically from the preceding paragraph.

# -- First new bullet: "Added notes to emphasise that the appropriate time period for comparison is one day."
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Text = "PLACEHOLDER"
$target = $d.Range($newRange.Start, $newRange.End)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve">Added notes to emphasise that the appropriate </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>time period</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t xml:space="preserve"> for comparison is one day.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml1)

# -- Second new bullet: "Added the option to switch between kilometres and miles for cycling and driving."
#    This one is a single plain run, so a direct Text assignment is enough
#    (InsertParagraphAfter already carried over the ListParagraph/numId 7/
#    Courier New formatting from the preceding paragraph).
$lastPara2 = $d.Paragraphs.Last
$lastPara2.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Last
$newPara2.Range.Text = "Added the option to switch between kilometres and miles for cycling and driving."

Write-Host "Edits applied"
